# Workbook reference
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. USMarket: reset its selection back to A1 (the default) while it is
#    still the active sheet, so that once it stops being the active tab
#    its tabSelected flag is dropped cleanly.
# ---------------------------------------------------------------------
$usSheet = $wb.Worksheets.Item(1)
[void]$usSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Create the "EuropeMarket" sheet right after USMarket
# ---------------------------------------------------------------------
$europeSheet = $wb.Worksheets.Add($null, $usSheet)
$europeSheet.Name = "EuropeMarket"

$europeSheet.Range("A1").Value = "Assertions"
$europeSheet.Range("A2").Value = "FTSE 100 Index"
$europeSheet.Range("A3").Value = "DAX"
$europeSheet.Range("A4").Value = "CAC 40 Index"
$europeSheet.Range("A5").Value = "FTSE MIB Index"
$europeSheet.Range("A6").Value = "IBEX 35 Index"
$europeSheet.Range("A7").Value = "STOXX Europe 600 Index"

$europeSheet.Columns.Item(1).ColumnWidth = 20.53
[void]$europeSheet.Range("A8").Select()
$europeSheet.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 3. Create the "AsiaMarket" sheet right after EuropeMarket
# ---------------------------------------------------------------------
$asiaSheet = $wb.Worksheets.Add($null, $europeSheet)
$asiaSheet.Name = "AsiaMarket"

$asiaSheet.Range("A1").Value = "Assertions"
$asiaSheet.Range("A2").Value = "The Asia Dow Index USD"
$asiaSheet.Range("A3").Value = "NIKKEI 225 Index"
$asiaSheet.Range("A4").Value = "Hang Seng Index"
$asiaSheet.Range("A5").Value = "Shanghai Composite Index"
$asiaSheet.Range("A6").Value = "S&P BSE Sensex Index"
$asiaSheet.Range("A7").Value = "FTSE Straits Times Index"

$asiaSheet.Columns.Item(1).ColumnWidth = 22.26
[void]$asiaSheet.Range("A8").Select()
$asiaSheet.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 4. Make AsiaMarket the active tab (this clears tabSelected on USMarket
#    and sets workbookView activeTab accordingly)
# ---------------------------------------------------------------------
[void]$asiaSheet.Activate()

Write-Host "Workbook now has $($wb.Worksheets.Count) sheets"
